$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

$ws.Range("C1").Value = "user1"
$ws.Range("E1").Value = "password1"
